$wb = $excel.ActiveWorkbook

# --- Insert the new worksheet right after "IL2_mRNA_new_data" ---
$anchor = $wb.Worksheets.Item("IL2_mRNA_new_data")
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "IL2_mRNA_new_data_20250725"
$anchor = $wb.Worksheets.Item("IL2_mRNA_new_data")
$newSheet.Move($null, $anchor)

# Re-fetch a fresh reference to the moved/renamed sheet (avoids stale handles).
$ws = $wb.Worksheets.Item("IL2_mRNA_new_data_20250725")

# --- Header row ---
$ws.Range("A1").Value = "condition"
$ws.Range("B1").Value = "treatment"
$ws.Range("C1").Value = "actin"
$ws.Range("D1").Value = "il2"
$ws.Range("E1").Value = "ratio"

# --- Data rows (condition, treatment, ratio) ---
$data = @(
    @("CM", "CM", 0.011723329425556858),
    @("CM", "CM", 0.01452081316553727),
    @("CM", "CM", 0.020387359836901122),
    @("CM", "CM", 0.017235),
    @("CM", "CM", 0.01579),
    @("PI", "PI", 1.4489795918367347),
    @("PI", "PI", 1.451911935110081),
    @("PI", "PI", 1.5771065182829889),
    @("PI", "PI", 1.468903),
    @("PI", "PI", 1.49724),
    @("PI", "PI_VCT", 1.2161458333333333),
    @("PI", "PI_VCT", 1.6234756097560976),
    @("PI", "PI_VCT", 1.6670918367346939),
    @("PI", "PI_VCT", 1.51834),
    @("PI", "PI_VCT", 1.498932),
    @("PI", "PI_100", 1.539235412474849),
    @("PI", "PI_100", 1.6887661141804788),
    @("PI", "PI_100", 0.835243553008596),
    @("PI", "PI_100", 1.34669),
    @("PI", "PI_100", 1.350205),
    @("PI", "PI_250", 0.5528455284552846),
    @("PI", "PI_250", 0.27899846704138986),
    @("PI", "PI_250", 0.47662018047579985),
    @("PI", "PI_250", 0.435023),
    @("PI", "PI_250", 0.445758),
    @("PI", "PI_500", 0.5227790432801822),
    @("PI", "PI_500", 0.691089108910891),
    @("PI", "PI_500", 0.48097826086956524),
    @("PI", "PI_500", 0.566432),
    @("PI", "PI_500", 0.559981),
    @("PI", "PI_1000", 0.2448224852071006),
    @("PI", "PI_1000", 0.18401332223147376),
    @("PI", "PI_1000", 0.4287003610108303),
    @("PI", "PI_1000", 0.32414),
    @("PI", "PI_1000", 0.278255),
    @("PI", "PI_1500", 0.11313591495823842),
    @("PI", "PI_1500", 0.1869530628480509),
    @("PI", "PI_1500", 0.13944817300521997),
    @("PI", "PI_1500", 0.14099),
    @("PI", "PI_1500", 0.139885)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $r = $r + 1
}

# --- View state for the previously-selected "IL2_mRNA_new_data" sheet ---
$old = $wb.Worksheets.Item("IL2_mRNA_new_data")
$old.Range("A1:B25").Select()

# --- View state for the new sheet: active cell D1, becomes the selected tab ---
$ws.Range("D1").Select()
$ws.Activate()
